$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 holding issue #12. Column A ("12") must stay text (matching the
# other rows, which store Issue ID as inline/shared text, not a number), so
# force the cell to text format before assigning the numeric-looking value.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "12"
$ws.Range("B8").Value = "newrelic_alert_channel Resource newrelic_synthetics_monitor: Cannot unset validation_string"
$ws.Range("C8").Value = "open"
$ws.Range("D8").Value = "2025-03-24T09:04:06Z"
$ws.Range("E8").Value = "bug"
